$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to remain text, matching the workbook's existing inline-string cells
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("B2").Value = 'Bitcoin'
$ws.Range("C2").Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range("D2").Value = '31.245.16'
$ws.Range("E2").Value = '  +2.24%  '
$ws.Range("B3").Value = 'Ethereum'
$ws.Range("C3").Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range("D3").Value = '1.998.55'
$ws.Range("E3").Value = '  +6.15%  '
$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("C4").Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("B5").Value = 'XRP'
$ws.Range("C5").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D5").Value = '0.7738'
$ws.Range("E5").Value = '  +63.52%  '
$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D6").Value = '254.53'
$ws.Range("E6").Value = '  +3.30%  '
$ws.Range("B7").Value = 'USDC'
$ws.Range("C7").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D7").Value = '1.000'
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D8").Value = '0.3490'
$ws.Range("E8").Value = '  +20.63%  '
$ws.Range("B9").Value = 'Solana'
$ws.Range("C9").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D9").Value = '27.76'
$ws.Range("E9").Value = '  +24.68%  '
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").Value = '0.07056'
$ws.Range("E10").Value = '  +7.93%  '
$ws.Range("B11").Value = 'Polygon'
$ws.Range("C11").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D11").Value = '0.8424'
$ws.Range("E11").Value = '  +10.74%  '
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").Value = '0.08211'
$ws.Range("E12").Value = '  +4.91%  '
$ws.Range("B13").Value = 'Litecoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D13").Value = '100.93'
$ws.Range("E13").Value = '  +1.20%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '5.648'
$ws.Range("E14").Value = '  +7.89%  '
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '1.996.39'
$ws.Range("E15").Value = '  +6.08%  '
$ws.Range("B16").Value = 'Avalanche'
$ws.Range("C16").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D16").Value = '15.28'
$ws.Range("E16").Value = '  +15.93%  '
$ws.Range("B17").Value = 'BitcoinCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D17").Value = '273.20'
$ws.Range("E17").Value = '  -3.88%  '
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '31.239.75'
$ws.Range("E18").Value = '  +2.30%  '
$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").Value = '5.968'
$ws.Range("E19").Value = '  +11.67%  '
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").Value = '0.000008019'
$ws.Range("E20").Value = '  +6.63%  '
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '2.256.26'
$ws.Range("E21").Value = '  +6.21%  '
$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").Value = '0.9999'
$ws.Range("E22").Value = '  +0.14%  '
$ws.Range("B23").Value = 'BinanceUSD'
$ws.Range("C23").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("B24").Value = 'Chainlink'
$ws.Range("C24").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D24").Value = '7.105'
$ws.Range("E24").Value = '  +10.59%  '
$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D25").Value = '9.997'
$ws.Range("E25").Value = '  +9.02%  '
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").Value = '165.93'
$ws.Range("E26").Value = '  +1.49%  '
$ws.Range("B27").Value = 'Stellar'
$ws.Range("C27").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D27").Value = '0.1430'
$ws.Range("E27").Value = '  +46.37%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = '19.94'
$ws.Range("E28").Value = '  +4.87%  '
$ws.Range("B29").Value = 'LidoDAOToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D29").Value = '2.409'
$ws.Range("E29").Value = '  +26.39%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").Value = '1.595'
$ws.Range("E30").Value = '  +6.23%  '
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").Value = '1.369'
$ws.Range("E31").Value = '  +3.06%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '4.606'
$ws.Range("E32").Value = '  +8.38%  '
$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").Value = '4.444'
$ws.Range("E33").Value = '  +6.19%  '
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").Value = '0.05276'
$ws.Range("E34").Value = '  +8.67%  '
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").Value = '0.7879'
$ws.Range("E35").Value = '  +12.72%  '
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").Value = '1.229'
$ws.Range("E36").Value = '  +8.67%  '
$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").Value = '2.768'
$ws.Range("E37").Value = '  +0.19%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.02002'
$ws.Range("E38").Value = '  +5.10%  '
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").Value = '2.920'
$ws.Range("E39").Value = '  +1.58%  '
$ws.Range("B40").Value = 'Aave'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D40").Value = '81.76'
$ws.Range("E40").Value = '  +8.54%  '
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").Value = '6.737'
$ws.Range("E41").Value = '  +6.77%  '
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").Value = '0.4679'
$ws.Range("E42").Value = '  +10.16%  '
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").Value = '2.139'
$ws.Range("E43").Value = '  +8.35%  '
$ws.Range("B44").Value = 'TrustWalletToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D44").Value = '0.8581'
$ws.Range("E44").Value = '  +2.50%  '
$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").Value = '104.98'
$ws.Range("E45").Value = '  +3.51%  '
$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").Value = '1.001'
$ws.Range("E46").Value = '  -0.04%  '
$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D47").Value = '7.715'
$ws.Range("E47").Value = '  +9.93%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '9.985'
$ws.Range("E48").Value = '  +0.38%  '
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").Value = '37.32'
$ws.Range("E49").Value = '  +5.89%  '
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").Value = '1.571'
$ws.Range("E50").Value = '  +17.34%  '
$ws.Range("B51").Value = 'Decentraland'
$ws.Range("C51").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D51").Value = '0.4325'
$ws.Range("E51").Value = '  +9.33%  '
